$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail";
# becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 duplicated the old "Contact" row entirely - remove it,
# shifting the remaining rows (Description, Purpose, ...) up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Root extension row: Short / Definition now describe this specific
# extension instead of the generic "Extension" / "An Extension" text.
$elements.Range("K2").Value = "Parent Organization Hierarchy Level Description"
$elements.Range("L2").Value = "Description of the level of the parent within the organinzational hierarchy"

# Column K (Short) widens to fit the new, longer text (best-fit ~45.33 chars)
$elements.Columns.Item(11).ColumnWidth = 44.5
